$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 24 (M07) first, while the description is still in column E
$ws.Range("A24").Value = "M07"
$ws.Range("B24").Value = "18K"
$ws.Range("C24").Value = 570
$ws.Range("D24").Value = 500
$ws.Range("E24").Value = "i_temp wird nun im Porgamm negativ berechnet "

# Move the existing "Beschreibung" column one column to the right
# (to F) to make room for the new "Verstärker" column. Done in pieces
# to avoid touching the empty row 25 gap between the data blocks.
$ws.Range("E17:E24").Cut($ws.Range("F17:F24"))
$ws.Range("E26").Cut($ws.Range("F26"))

# New column header
$ws.Range("E17").Value = "Verstärker"

# Fill "ohne" for all existing measurement rows in the new column
$ws.Range("E18").Value = "ohne"
$ws.Range("E19").Value = "ohne"
$ws.Range("E20").Value = "ohne"
$ws.Range("E21").Value = "ohne"
$ws.Range("E22").Value = "ohne"
$ws.Range("E23").Value = "ohne"
$ws.Range("E24").Value = "ohne"
$ws.Range("E26").Value = "ohne"

# Update description text for D01 (row 26) to reference M06
$ws.Range("F26").Value = "Messung und Berechnung durch den µC, Übertragung der Leistungsdaten (Ref: M06)"

# Add new row 27: D02
$ws.Range("A27").Value = "D02"
$ws.Range("B27").Value = "18K"
$ws.Range("C27").Value = 570
$ws.Range("D27").Value = 500
$ws.Range("E27").Value = "ohne"
$ws.Range("F27").Value = "Berechnete Leisuntgsdaten mit vorher gedrehtem Vorzeichen von i_temp (Ref: M07)"

# Update selection to match target state
$ws.Range("I30").Select()
